$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -12.4764
$ws.Range("A3").Value = -22.0452
$ws.Range("A14").Value = -21.7468
$ws.Range("A16").Value = -21.75709999999999
$ws.Range("C18").Value = -11.7081
$ws.Range("A21").Value = -20.15689999999998
$ws.Range("A23").Value = -20.27779999999997
$ws.Range("C24").Value = -13.23569999999999
$ws.Range("A25").Value = -21.66249999999999
$ws.Range("C25").Value = -12.7937
$ws.Range("A26").Value = -21.17429999999997
$ws.Range("C27").Value = -12.6068
$ws.Range("A29").Value = -20.86649999999999
$ws.Range("C30").Value = -12.9026
$ws.Range("C31").Value = -13.3699
$ws.Range("C39").Value = -12.4937
$ws.Range("A40").Value = -19.96549999999999
$ws.Range("C42").Value = -12.9179
$ws.Range("C48").Value = -11.3808
$ws.Range("C51").Value = -11.4435
$ws.Range("C52").Value = -11.26169999999999
$ws.Range("A53").Value = -21.7958
$ws.Range("C55").Value = -13.7358
$ws.Range("C56").Value = -12.28779999999999
$ws.Range("A57").Value = -21.8965
$ws.Range("C57").Value = -12.57159999999999
$ws.Range("A59").Value = -22.3753
$ws.Range("C60").Value = -13.2664
$ws.Range("A65").Value = -21.87939999999998
$ws.Range("A69").Value = -21.5981
$ws.Range("C73").Value = -12.8025
$ws.Range("C74").Value = -12.60570000000001
$ws.Range("A79").Value = -20.62040000000001
$ws.Range("A83").Value = -21.72519999999999
$ws.Range("C89").Value = -10.56720000000001
$ws.Range("C90").Value = -12.4145
$ws.Range("A91").Value = -21.45090000000002
$ws.Range("C92").Value = -10.89599999999999
$ws.Range("A93").Value = -21.12119999999999
$ws.Range("A100").Value = -21.80419999999999
